# Weekly update: Fruta / hortaliza -- Vega Monumental Concepción, Caqui
# Applies the row-content changes described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44707
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 12500
$ws.Range("Q2").Value = '$/caja 12 kilos empedrada'
$ws.Range("R2").Value = 'Provincia de Curicó'
$ws.Range("S2").Value = 1042
$ws.Range("T2").Value = 12
# Row 3
$ws.Range("D3").Value = 45100
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 18
# Row 4
$ws.Range("D4").Value = 44330
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 15500
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("R4").Value = 'Provincia de Curicó'
$ws.Range("S4").Value = 861
$ws.Range("T4").Value = 18
# Row 5
$ws.Range("D5").Value = 45077
$ws.Range("M5").Value = 140
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 12857
$ws.Range("Q5").Value = '$/caja 12 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 12857
$ws.Range("T5").Value = 1
# Row 6
$ws.Range("D6").Value = 45077
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 11000
$ws.Range("P6").Value = 11000
$ws.Range("Q6").Value = '$/caja 12 kilos granel'
$ws.Range("S6").Value = 11000
$ws.Range("T6").Value = 1
# Row 7
$ws.Range("D7").Value = 45084
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 17000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 17500
$ws.Range("S7").Value = 972
# Row 10
$ws.Range("D10").Value = 45093
$ws.Range("M10").Value = 140
$ws.Range("N10").Value = 17000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 17429
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("S10").Value = 968
$ws.Range("T10").Value = 18
# Row 11
$ws.Range("D11").Value = 45092
$ws.Range("M11").Value = 140
$ws.Range("N11").Value = 18000
$ws.Range("O11").Value = 19000
$ws.Range("P11").Value = 18429
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("R11").Value = 'Provincia de Curicó'
$ws.Range("S11").Value = 1024
$ws.Range("T11").Value = 18
# Row 12
$ws.Range("D12").Value = 44708
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 12571
$ws.Range("Q12").Value = '$/caja 12 kilos empedrada'
$ws.Range("S12").Value = 1048
$ws.Range("T12").Value = 12
# Row 13
$ws.Range("D13").Value = 44334
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 11000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 11500
$ws.Range("Q13").Value = '$/caja 12 kilos granel'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 11500
$ws.Range("T13").Value = 1
# Row 14
$ws.Range("D14").Value = 45097
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 19000
$ws.Range("R14").Value = 'Región del Maule'
$ws.Range("S14").Value = 1056
# Row 17
$ws.Range("D17").Value = 44719
$ws.Range("M17").Value = 50
$ws.Range("P17").Value = 14400
$ws.Range("R17").Value = 'Región del Maule'
$ws.Range("S17").Value = 800
# Row 18
$ws.Range("D18").Value = 44742
$ws.Range("L18").Value = 'Segunda'
$ws.Range("N18").Value = 14000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 14500
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 806
# Row 19
$ws.Range("D19").Value = 44714
$ws.Range("L19").Value = 'Primera'
